# Rename the "form_id" setting to "table_id" in the settings sheet so that
# the generated definitions.csv / properties.csv use the table_id, and make
# the settings sheet the active/selected tab (matching the saved view state
# of the authored workbook).

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")

# "form_id" -> "table_id" (keep the value "breathcounter" in B2 untouched).
$wsSettings.Range("A2").Value = "table_id"

# The renamed setting is a "group header" row like A1 (setting_name) and
# A4 (survey), so pick up that header formatting instead of the plain
# A3 (form_version) formatting it inherited while it was "form_id".
$wsSettings.Range("A2").Style = $wsSettings.Range("A4").Style

# The workbook was re-saved with the "settings" sheet active (instead of
# "survey"), with the selection left on A3.
$wsSettings.Activate()
$wsSettings.Range("A3").Select()
